# Auto-generated edit script: update "想去人数" (F column) counts
# across sheets 展览, 演出, 本地生活, 全部类型, matching the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 27079
$ws.Range("F4").Value = 668
$ws.Range("F5").Value = 201
$ws.Range("F6").Value = 583
$ws.Range("F8").Value = 385
$ws.Range("F9").Value = 499
$ws.Range("F11").Value = 55
$ws.Range("F12").Value = 320
$ws.Range("F13").Value = 105
$ws.Range("F14").Value = 523
$ws.Range("F16").Value = 1660
$ws.Range("F17").Value = 275
$ws.Range("F18").Value = 1088
$ws.Range("F19").Value = 201
$ws.Range("F22").Value = 112

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 461
$ws.Range("F17").Value = 30

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 283

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 283
$ws.Range("F5").Value = 27079
$ws.Range("F7").Value = 668
$ws.Range("F10").Value = 201
$ws.Range("F16").Value = 461
$ws.Range("F17").Value = 583
$ws.Range("F21").Value = 385
$ws.Range("F22").Value = 499
$ws.Range("F24").Value = 55
$ws.Range("F26").Value = 320
$ws.Range("F27").Value = 105
$ws.Range("F30").Value = 523
$ws.Range("F33").Value = 1660
$ws.Range("F34").Value = 275
$ws.Range("F35").Value = 1088
$ws.Range("F36").Value = 30
$ws.Range("F37").Value = 201
$ws.Range("F40").Value = 112
